# Fruta / hortaliza, semanal
# Re-order the weekly price records (rows 2-12) into their new positions.
# Columns A, B, C, E, F, G, I, R are identical across all these rows and
# are left untouched; only D, H, J, K, L, M, N, O, P, Q need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for the columns that vary, keyed by row.
$data = @{
    2  = @{ D = 44701; H = "Española"; J = 400; K = 19000; L = 20000; M = 19500; N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 650; Q = 30 }
    3  = @{ D = 44438; H = "Española"; J = 400; K = 11000; L = 12000; M = 11500; N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 383; Q = 30 }
    4  = @{ D = 44427; H = "Madrigal"; J = 400; K = 12000; L = 13000; M = 12500; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 312; Q = 40 }
    5  = @{ D = 44426; H = "Española"; J = 600; K = 11500; L = 12000; M = 11750; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 392; Q = 30 }
    6  = @{ D = 44426; H = "Madrigal"; J = 500; K = 12500; L = 13000; M = 12750; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 319; Q = 40 }
    7  = @{ D = 44484; H = "Española"; J = 300; K = 9000;  L = 10000; M = 9500;  N = "`$/caja 30 unidades"; O = "Provincia del Elquí"; P = 317; Q = 30 }
    8  = @{ D = 44729; H = "Madrigal"; J = 400; K = 16000; L = 17000; M = 16500; N = "`$/caja 40 unidades"; O = "Provincia del Elquí"; P = 412; Q = 40 }
    9  = @{ D = 44498; H = "Española"; J = 400; K = 8500;  L = 9000;  M = 8750;  N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 292; Q = 30 }
    10 = @{ D = 44420; H = "Madrigal"; J = 800; K = 14000; L = 15000; M = 14500; N = "`$/caja 40 unidades"; O = "Provincia de Limarí"; P = 362; Q = 40 }
    11 = @{ D = 44420; H = "Madrigal"; J = 700; K = 13000; L = 14000; M = 13500; N = "`$/caja 40 unidades"; O = "Provincia del Elquí"; P = 338; Q = 40 }
    12 = @{ D = 44687; H = "Española"; J = 400; K = 18000; L = 19000; M = 18500; N = "`$/caja 30 unidades"; O = "Provincia de Limarí"; P = 617; Q = 30 }
}

# Mapping of original row -> new row (a permutation of rows 2..12).
$mapping = @{
    2  = 5
    3  = 6
    4  = 9
    5  = 11
    6  = 12
    7  = 10
    8  = 8
    9  = 7
    10 = 3
    11 = 4
    12 = 2
}

foreach ($srcRow in $mapping.Keys) {
    $dstRow = $mapping[$srcRow]
    $rec = $data[$srcRow]

    $ws.Range("D$dstRow").Value = $rec.D
    $ws.Range("H$dstRow").Value = $rec.H
    $ws.Range("J$dstRow").Value = $rec.J
    $ws.Range("K$dstRow").Value = $rec.K
    $ws.Range("L$dstRow").Value = $rec.L
    $ws.Range("M$dstRow").Value = $rec.M
    $ws.Range("N$dstRow").Value = $rec.N
    $ws.Range("O$dstRow").Value = $rec.O
    $ws.Range("P$dstRow").Value = $rec.P
    $ws.Range("Q$dstRow").Value = $rec.Q
}
